$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$ws = $wb.Worksheets.Item("BDSBaPCF")

# "biomass" no longer bids at its peak capacity factor - flip the flag off.
# (B17 is "=B9", so it recalculates to 0 automatically.)
$ws.Range("B9").Value = 0

# Update the sheet's remembered selection (B5 -> B10), then restore the
# originally active sheet ("About") so the workbook's active tab is unchanged.
$ws.Activate() | Out-Null
$ws.Range("B10").Select() | Out-Null
$wsAbout.Activate() | Out-Null
